# Update "想去人数" (number of people interested) counts in column F
# for the affected rows across the "展览", "演出" and "全部类型" sheets.
# (The "本地生活" sheet is untouched - it has no data rows.)

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 2741
$ws1.Range("F8").Value  = 1665
$ws1.Range("F9").Value  = 1743
$ws1.Range("F11").Value = 280
$ws1.Range("F12").Value = 710
$ws1.Range("F13").Value = 871
$ws1.Range("F14").Value = 149
$ws1.Range("F19").Value = 534
$ws1.Range("F20").Value = 6277
$ws1.Range("F21").Value = 244
$ws1.Range("F22").Value = 1351
$ws1.Range("F23").Value = 138
$ws1.Range("F26").Value = 294
$ws1.Range("F27").Value = 249
$ws1.Range("F29").Value = 1089
$ws1.Range("F30").Value = 885
$ws1.Range("F32").Value = 82
$ws1.Range("F34").Value = 453
$ws1.Range("F35").Value = 1311
$ws1.Range("F37").Value = 137
$ws1.Range("F38").Value = 210
$ws1.Range("F41").Value = 172
$ws1.Range("F42").Value = 142

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 2

# --- Sheet "全部类型" (All types, combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 2741
$ws4.Range("F11").Value = 1665
$ws4.Range("F12").Value = 1743
$ws4.Range("F14").Value = 280
$ws4.Range("F15").Value = 710
$ws4.Range("F17").Value = 871
$ws4.Range("F18").Value = 149
$ws4.Range("F22").Value = 534
$ws4.Range("F23").Value = 6277
$ws4.Range("F24").Value = 244
$ws4.Range("F25").Value = 1351
$ws4.Range("F26").Value = 2
$ws4.Range("F27").Value = 138
$ws4.Range("F30").Value = 294
$ws4.Range("F31").Value = 249
$ws4.Range("F33").Value = 1089
$ws4.Range("F34").Value = 885
$ws4.Range("F36").Value = 82
$ws4.Range("F38").Value = 453
$ws4.Range("F39").Value = 1311
$ws4.Range("F41").Value = 137
$ws4.Range("F42").Value = 210
$ws4.Range("F45").Value = 172
$ws4.Range("F49").Value = 142
